$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Fri Oct 06 11:32:43 EDT 2023"
$ws.Range("B3").Value = "Fri Oct 06 11:32:57 EDT 2023"
$ws.Range("B4").Value = "Fri Oct 06 11:33:10 EDT 2023"
